{"js": "// Update the 25 division-fact answers in the single 5-column practice\n// table. The table is laid out as 5 \"problem\" rows (table row indices\n// 0, 4, 8, 12, 16), each followed by 3 blank rows, with 5 cells per\n// problem row (columns 0-4). Each entry below gives the zero-based\n// table row/column together with the expected old text (for a sanity\n// check) and the new text to write.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"89\u00f72=44, 1\", newText: \"48\u00f79=5, 3\" },\n  { row: 0, col: 1, oldText: \"60\u00f76=10, 0\", newText: \"66\u00f74=16, 2\" },\n  { row: 0, col: 2, oldText: \"33\u00f79=3, 6\", newText: \"86\u00f75=17, 1\" },\n  { row: 0, col: 3, oldText: \"44\u00f77=6, 2\", newText: \"30\u00f78=3, 6\" },\n  { row: 0, col: 4, oldText: \"86\u00f72=43, 0\", newText: \"93\u00f73=31, 0\" },\n\n  { row: 4, col: 0, oldText: \"67\u00f73=22, 1\", newText: \"11\u00f75=2, 1\" },\n  { row: 4, col: 1, oldText: \"20\u00f74=5, 0\", newText: \"11\u00f76=1, 5\" },\n  { row: 4, col: 2, oldText: \"98\u00f76=16, 2\", newText: \"84\u00f76=14, 0\" },\n  { row: 4, col: 3, oldText: \"35\u00f77=5, 0\", newText: \"71\u00f76=11, 5\" },\n  { row: 4, col: 4, oldText: \"91\u00f74=22, 3\", newText: \"55\u00f73=18, 1\" },\n\n  { row: 8, col: 0, oldText: \"21\u00f77=3, 0\", newText: \"37\u00f74=9, 1\" },\n  { row: 8, col: 1, oldText: \"89\u00f77=12, 5\", newText: \"53\u00f75=10, 3\" },\n  { row: 8, col: 2, oldText: \"18\u00f73=6, 0\", newText: \"13\u00f77=1, 6\" },\n  { row: 8, col: 3, oldText: \"13\u00f72=6, 1\", newText: \"50\u00f77=7, 1\" },\n  { row: 8, col: 4, oldText: \"52\u00f76=8, 4\", newText: \"58\u00f78=7, 2\" },\n\n  { row: 12, col: 0, oldText: \"33\u00f76=5, 3\", newText: \"45\u00f76=7, 3\" },\n  { row: 12, col: 1, oldText: \"18\u00f74=4, 2\", newText: \"44\u00f76=7, 2\" },\n  { row: 12, col: 2, oldText: \"41\u00f78=5, 1\", newText: \"52\u00f76=8, 4\" },\n  { row: 12, col: 3, oldText: \"42\u00f74=10, 2\", newText: \"21\u00f78=2, 5\" },\n  { row: 12, col: 4, oldText: \"22\u00f72=11, 0\", newText: \"54\u00f77=7, 5\" },\n\n  { row: 16, col: 0, oldText: \"29\u00f78=3, 5\", newText: \"19\u00f74=4, 3\" },\n  { row: 16, col: 1, oldText: \"73\u00f77=10, 3\", newText: \"67\u00f78=8, 3\" },\n  { row: 16, col: 2, oldText: \"15\u00f78=1, 7\", newText: \"82\u00f78=10, 2\" },\n  { row: 16, col: 3, oldText: \"54\u00f72=27, 0\", newText: \"95\u00f73=31, 2\" },\n  { row: 16, col: 4, oldText: \"22\u00f72=11, 0\", newText: \"46\u00f73=15, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load the current value of every target cell up front so we can\n// validate before writing (defends against an unexpected layout).\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { oldText, newText } = replacements[i];\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${replacements[i].row}, col ${replacements[i].col}: ` +\n        `expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 division-fact answers in the single 5-column practice\n# table. The table has 20 rows: 5 \"problem\" rows (1-based table row\n# indices 1, 5, 9, 13, 17), each followed by 3 blank rows, with 5 cells\n# per problem row (columns 1-5). Each entry below gives the 1-based\n# table row/column together with the expected old text (for a sanity\n# check) and the new text to write.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"89\u00f72=44, 1\";  New = \"48\u00f79=5, 3\" },\n    @{ Row = 1;  Col = 2; Old = \"60\u00f76=10, 0\";  New = \"66\u00f74=16, 2\" },\n    @{ Row = 1;  Col = 3; Old = \"33\u00f79=3, 6\";   New = \"86\u00f75=17, 1\" },\n    @{ Row = 1;  Col = 4; Old = \"44\u00f77=6, 2\";   New = \"30\u00f78=3, 6\" },\n    @{ Row = 1;  Col = 5; Old = \"86\u00f72=43, 0\";  New = \"93\u00f73=31, 0\" },\n\n    @{ Row = 5;  Col = 1; Old = \"67\u00f73=22, 1\";  New = \"11\u00f75=2, 1\" },\n    @{ Row = 5;  Col = 2; Old = \"20\u00f74=5, 0\";   New = \"11\u00f76=1, 5\" },\n    @{ Row = 5;  Col = 3; Old = \"98\u00f76=16, 2\";  New = \"84\u00f76=14, 0\" },\n    @{ Row = 5;  Col = 4; Old = \"35\u00f77=5, 0\";   New = \"71\u00f76=11, 5\" },\n    @{ Row = 5;  Col = 5; Old = \"91\u00f74=22, 3\";  New = \"55\u00f73=18, 1\" },\n\n    @{ Row = 9;  Col = 1; Old = \"21\u00f77=3, 0\";   New = \"37\u00f74=9, 1\" },\n    @{ Row = 9;  Col = 2; Old = \"89\u00f77=12, 5\";  New = \"53\u00f75=10, 3\" },\n    @{ Row = 9;  Col = 3; Old = \"18\u00f73=6, 0\";   New = \"13\u00f77=1, 6\" },\n    @{ Row = 9;  Col = 4; Old = \"13\u00f72=6, 1\";   New = \"50\u00f77=7, 1\" },\n    @{ Row = 9;  Col = 5; Old = \"52\u00f76=8, 4\";   New = \"58\u00f78=7, 2\" },\n\n    @{ Row = 13; Col = 1; Old = \"33\u00f76=5, 3\";   New = \"45\u00f76=7, 3\" },\n    @{ Row = 13; Col = 2; Old = \"18\u00f74=4, 2\";   New = \"44\u00f76=7, 2\" },\n    @{ Row = 13; Col = 3; Old = \"41\u00f78=5, 1\";   New = \"52\u00f76=8, 4\" },\n    @{ Row = 13; Col = 4; Old = \"42\u00f74=10, 2\";  New = \"21\u00f78=2, 5\" },\n    @{ Row = 13; Col = 5; Old = \"22\u00f72=11, 0\";  New = \"54\u00f77=7, 5\" },\n\n    @{ Row = 17; Col = 1; Old = \"29\u00f78=3, 5\";   New = \"19\u00f74=4, 3\" },\n    @{ Row = 17; Col = 2; Old = \"73\u00f77=10, 3\";  New = \"67\u00f78=8, 3\" },\n    @{ Row = 17; Col = 3; Old = \"15\u00f78=1, 7\";   New = \"82\u00f78=10, 2\" },\n    @{ Row = 17; Col = 4; Old = \"54\u00f72=27, 0\";  New = \"95\u00f73=31, 2\" },\n    @{ Row = 17; Col = 5; Old = \"22\u00f72=11, 0\";  New = \"46\u00f73=15, 1\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.Old) {\n        throw \"Unexpected cell text at row $($r.Row), col $($r.Col): expected '$($r.Old)' but found '$current'\"\n    }\n    $cell.Range.Text = $r.New\n}\n"}
